# Generate Report for Handoff
# The "b.md" row is now ready for handoff: status flips to "Ready for handoff",
# a fresh handoff xliff was generated for both locales, and the zh-cn locale's
# handback turned out to be stale against latest source, so an error detail is
# recorded and its "Content Duplicate" flag clears.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$genDate = "2016-08-15 20:32:14"

$zhHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate = "2016-08-15 20:32:06"
$deHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate = $genDate

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4629a1e6c2df686dfd9b364929b7c7602f3e3b3/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d590bb729015c510e62d01b0f65b02d42e923dd1/e2e/b.md."

# --- Overview sheet: row 3 is "b.md" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $genDate

# --- zh-cn sheet: row 3 is "b.md" ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $status
# Plain "True"/"False" values auto-coerce to booleans; force text (as the
# source data uses) via an apostrophe-prefixed literal, then restore the
# default style so no stray formatting/number-format is left behind.
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("F3").Style = "Normal"
$wsZh.Range("G3").Value = $zhHandoffFile
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = 39.1666666666667

# --- de-de sheet: row 3 is "b.md" ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $status
$wsDe.Range("G3").Value = $deHandoffFile
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.1666666666667
